# Rebuild the "registros" table: row 2 (Item=1, aciclovir/250mg/ml/cimed/143810181/Pendente)
# is replaced by four rows (Item 6..9) describing CLORETO DE SODIO;GLICOSE / SORO GLICOSADO,
# all supplied by JP INDUSTRIA FARMACEUTICA S/A, still Pendente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Force a digit-only string (e.g. a registry number) to be stored as TEXT
    # instead of being auto-coerced to a number, then drop back to the default
    # (unstyled) cell style so no stray number format sticks around.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$rows = @(
    @{ Row = 2; Item = 6; Desc = "CLORETO DE SÓDIO;GLICOSE"; Conc = "1000ml"; Registro = "104910019" },
    @{ Row = 3; Item = 7; Desc = "CLORETO DE SÓDIO;GLICOSE"; Conc = "250ml";  Registro = "104910019" },
    @{ Row = 4; Item = 8; Desc = "CLORETO DE SÓDIO;GLICOSE"; Conc = "500ml";  Registro = "104910019" },
    @{ Row = 5; Item = 9; Desc = "SORO GLICOSADO 5%, BOLSA SISTEMA FECHADO, FRASCO COM 1000ML"; Conc = "1000ml"; Registro = "104910020" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.Desc
}
foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("C$n").Value = $r.Conc
}
foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("D$n").Value = "JP INDUSTRIA FARMACEUTICA S/A"
}
foreach ($r in $rows) {
    $n = $r.Row
    Set-TextValue $ws.Range("E$n") $r.Registro
}
foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("F$n").Value = "Pendente"
}
foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = $r.Item
}
